# Capacitación a Docentes.pptx -- "Migración Docente, estudiante y usuario"
#
# Real content edits recovered from the target diff:
#   1. The "Provisionalmente usaremos" IP address on slide 10 changes
#      from 192.168.123.15/idiomas to 192.168.123.11/idiomas.
#   2. The cached text of every datetimeFigureOut date field (slide
#      master + all 11 slide layouts) is bumped from 10/06/2018 to
#      11/06/2018 (PowerPoint re-stamped these the next time the
#      deck was saved).
#
# (The remaining hunks in the diff -- the empty p15:sldGuideLst ext
# block on <p:presentation>, and the mc:AlternateContent namespace
# shuffling on every slide's transition -- are PowerPoint's own XML
# writer doing incidental, content-free re-serialization on save; they
# are not reachable through the documented PowerPoint object model, so
# they are intentionally left alone here.)

$p = $ppt.ActivePresentation

# --- 1. Fix the lab IP address referenced on the "Acceso al sistema" slide ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi, 1)
                if ($para.Text -eq "192.168.123.15/idiomas") {
                    $run = $para.Runs(1, 1)
                    $run.Text = "192.168.123.11/idiomas"
                }
            }
        }
    }
}

# --- 2. Bump the cached date field text on the master + every layout ---
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "10/06/2018") {
            $tr.Text = "11/06/2018"
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "10/06/2018") {
                $tr.Text = "11/06/2018"
            }
        }
    }
}
